# "Entregas na sprint 1"
# Mark additional "Sprint 2" deliveries as Ok on the "Presenças" sheet,
# and fill in the "Respostas_Entregues" sheet's second exercise column (C)
# to mirror the actual delivery status recorded for each person.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Presenças" ------------------------------------------------
$wsPresencas = $wb.Worksheets.Item(1)

$presencasRows = @(5, 11, 17, 18, 22, 23)
foreach ($r in $presencasRows) {
    $wsPresencas.Range("E$r").Value = "Ok"
}

# --- Sheet 2: "Respostas_Entregues" --------------------------------------
$wsRespostas = $wb.Worksheets.Item(2)

$respostasValues = @{
    2  = "Ok"
    3  = "F"
    4  = "Ok"
    5  = "Ok"
    6  = "F"
    7  = "F"
    8  = "F"
    9  = "Ok"
    10 = "J"
    11 = "F"
    12 = "F"
    13 = "Ok"
    14 = "F"
    15 = "J"
    16 = "Ok"
    17 = "F"
    18 = "F"
    19 = "F"
    20 = "Ok"
    21 = "F"
    22 = "F"
    23 = "F"
    24 = "Ok"
    25 = "Ok"
}

foreach ($r in 2..25) {
    $wsRespostas.Range("C$r").Value = $respostasValues[$r]
}

# --- Restore the active-cell selections recorded in the saved file -------
[void]$wsPresencas.Range("E21").Select()
[void]$wsRespostas.Range("C4").Select()
[void]$wsRespostas.Activate()
